# Auto-generated update script for cryptos.xlsx
# Applies the per-row Price (D) / Volume(1h) (E) updates and the two
# Coin/Link (B/C) row swaps described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to store the value as literal TEXT (not an
    # auto-coerced number/date) while leaving the cells style
    # untouched, matching the inline-string cells already in the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '67.948.35'
$ws.Cells.Item(2, 5).Value = '  +1.41%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '3.258.63'
$ws.Cells.Item(3, 5).Value = '  -0.28%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '584.45'
$ws.Cells.Item(5, 5).Value = '  +0.70%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '184.07'
$ws.Cells.Item(6, 5).Value = '  +3.69%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -1.57%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +3.70%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.61%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.416'
$ws.Cells.Item(11, 5).Value = '  +1.52%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '3.824.84'
$ws.Cells.Item(12, 5).Value = '  -0.19%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +0.38%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '28.54'
$ws.Cells.Item(14, 5).Value = '  +1.48%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '67.978.93'
$ws.Cells.Item(15, 5).Value = '  +1.50%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +2.72%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '3.250.12'
$ws.Cells.Item(17, 5).Value = '  -0.41%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '5.86'
$ws.Cells.Item(18, 5).Value = '  +0.43%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '13.62'
$ws.Cells.Item(19, 5).Value = '  +1.55%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '381.90'
$ws.Cells.Item(20, 5).Value = '  +3.15%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '7.68'
$ws.Cells.Item(21, 5).Value = '  +0.58%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.04%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '71.27'
$ws.Cells.Item(23, 5).Value = '  +0.57%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '0.515'
$ws.Cells.Item(24, 5).Value = '  +0.80%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.81%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.59%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +2.26%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '1.00'
$ws.Cells.Item(28, 5).Value = '  +0.15%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.29%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '5.69'
$ws.Cells.Item(30, 5).Value = '  +1.02%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '7.34'
$ws.Cells.Item(31, 5).Value = '  +8.18%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '22.89'
$ws.Cells.Item(32, 5).Value = '  +1.30%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -0.01%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '1.27'
$ws.Cells.Item(34, 5).Value = '  +1.56%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +3.55%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '162.49'
$ws.Cells.Item(36, 5).Value = '  -4.03%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.39%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.835'
$ws.Cells.Item(38, 5).Value = '  -2.86%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '26.66'
$ws.Cells.Item(39, 5).Value = '  -1.98%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '6.75'
$ws.Cells.Item(40, 5).Value = '  +5.06%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '4.61'
$ws.Cells.Item(41, 5).Value = '  +6.17%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +1.25%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '41.42'
$ws.Cells.Item(43, 5).Value = '  +2.19%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '25.44'
$ws.Cells.Item(44, 5).Value = '  +2.67%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Bittensor'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Cells.Item(45, 4) '346.22'
$ws.Cells.Item(45, 5).Value = '  +0.94%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Hedera'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Cells.Item(46, 4) '0.0687'
$ws.Cells.Item(46, 5).Value = '  +2.23%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '2.645.64'
$ws.Cells.Item(47, 5).Value = '  -4.33%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +1.64%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -1.31%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Arweave'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Cells.Item(50, 4) '31.68'
$ws.Cells.Item(50, 5).Value = '  +4.25%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'ONDO'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Cells.Item(51, 4) '0.995'
$ws.Cells.Item(51, 5).Value = '  +1.34%  '

Write-Output "Applied 50 row updates."
